# New datas with more features.
# Adds a new training-log row (row 71) to Sheet1, mirroring the formatting
# of the previous row (row 70), and updates the sheet view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting of the last existing data row (70) onto the new row (71)
# so the new row picks up the same cell styles (date format, wrap text, etc.)
$ws.Range("A70:L70").Copy()
$ws.Range("A71:L71").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Match the custom row height used by the surrounding rows.
$ws.Rows.Item(71).RowHeight = 82.5

# Fill in the new row's data.
$ws.Range("A71").Value = 43214.347222222219
$ws.Range("B71").Value = "分类"
$ws.Range("C71").Value = "14分类"
$ws.Range("D71").Value = "batch_size=100 low_nums=2 use_biases=yes use_bn_low=True dropout_low=0.7"
$ws.Range("E71").Value = "最高标签，重新训练，原始数据加上序列train-hjxh365-2018-4-16-day-high-sequence"
$ws.Range("F71").Value = 0.7
$ws.Range("G71").Value = 0.6
$ws.Range("H71").Value = 1
$ws.Range("I71").Value = 1
$ws.Range("K71").Value = "python feed_run.py --output_mode=classes --output_nodes=14 --input_nums=129 --input_nodes=129 --low_nums=2 --low_nodes=129 --low_fun=elu --use_bn_input=True --one_hot=True --input_fun=tanh --batch_size=100 --learning_rate=0.001 --train_mode=Adadelta --eval_size=5400 --test_size=1339 --use_biases=yes --dropout_low=0.7"
$ws.Range("L71").Value = "logs-hjxh-2018-4-24-class14-sequence-percent60"
$ws.Range("J71").Value = "经过约17小时，拟合精度还能提高，泛化精度无法提高了，还有点下降，两者已明显分化。看来加入单纯序列不利于泛化。"

# Update the active selection to match the new edit location.
[void]$ws.Range("J70").Select()
